$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G7").Value = "Data wrangling with dplyr"
$ws.Range("G9").Value = "Exercise: Wrangling the Amniote Life History Database"
$ws.Range("G10").Value = "Combining data sets"
$ws.Range("G11").Value = "Exercise: Temperature effects on egg laying dates"
$ws.Range("G12").Value = "Data visualisation with ggplot2"

$ws.Range("G12").Select()
